$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Add two new test case rows (105 and 106) for new search test cases.
# Values are written in the specific order below so that new entries land
# in the shared strings table in the same order as the target workbook.
$ws.Range("A105").Value = "TestCase_B104"
$ws.Range("A106").Value = "TestCase_B105"

$ws.Range("B105").Value = "OPQA-554"
$ws.Range("B106").Value = "OPQA-555"

$ws.Range("C106").Value = "Verify that record view page of a post gets displayed when user clicks on article title in POSTs search results page"
$ws.Range("C105").Value = "Verify that record view page of a post gets displayed when user clicks on article title in ALL  search results page"

$ws.Range("D105").Value = "Y"
$ws.Range("D106").Value = "Y"

# Copy formatting from row 104 (an existing similar row) to the new rows
$ws.Range("A104:E104").Copy()
$ws.Range("A105:E105").PasteSpecial(-4122)
$ws.Range("A106:E106").PasteSpecial(-4122)

# Update the active selection to match the target state
$ws.Range("C97").Select()
